# Append: 2025-11-22 18:23 JST
# Update the "取得日時" (retrieved timestamp) column A for every data row
# on the ランサーズ sheet from 2025-11-22 12:41:03 to 2025-11-22 18:23:09.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-11-22 12:41:03"
$newTimestamp = "2025-11-22 18:23:09"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
